$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("steel_prim")
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 79.4507036657
$ws.Range("B15").Value = 0
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 0
$ws.Range("B32").Value = 0

$ws = $wb.Worksheets.Item("paper")
$ws.Range("B2").Value = 2330.5420294995
$ws.Range("B3").Value = 311.8444718212
$ws.Range("B4").Value = 931.6389662662
$ws.Range("B5").Value = 495.6006341929
$ws.Range("B6").Value = 24442.1652085912
$ws.Range("B7").Value = 16.0965894055
$ws.Range("B8").Value = 622.299137233
$ws.Range("B9").Value = 6932.5421345551
$ws.Range("B10").Value = 9970.731257272901
$ws.Range("B11").Value = 453.350733722
$ws.Range("B12").Value = 10332.3062555283
$ws.Range("B13").Value = 36.3951041709
$ws.Range("B14").Value = 7.0730814158
$ws.Range("B15").Value = 616.2000532224
$ws.Range("B16").Value = 3491.9743096676
$ws.Range("B17").Value = 5937.2597674029
$ws.Range("B18").Value = 3641.3479146008
$ws.Range("B19").Value = 1997.8478353432
$ws.Range("B20").Value = 392.8070473097
$ws.Range("B21").Value = 815.9434664426
$ws.Range("B22").Value = 859.972268227
$ws.Range("B23").Value = 13933.1439886563
$ws.Range("B24").Value = 13544.1149952597
$ws.Range("B25").Value = 5035.3691960091
$ws.Range("B26").Value = 2029.1459140583
$ws.Range("B27").Value = 1845.4135431481
$ws.Range("B29").Value = 29.3719505427
$ws.Range("B31").Value = 320.3985230944
$ws.Range("B33").Value = 3.426623842
$ws.Range("B34").Value = 110.3220201591
$ws.Range("B35").Value = 73.3031003756

$ws = $wb.Worksheets.Item("cement")
$ws.Range("B2").Value = 7963.7780625784
$ws.Range("B3").Value = 2480.4423465695
$ws.Range("B4").Value = 4401.2399948629
$ws.Range("B5").Value = 1967.470891339
$ws.Range("B6").Value = 30817.6245281912
$ws.Range("B7").Value = 3221.3292818066
$ws.Range("B8").Value = 13857.777530483
$ws.Range("B9").Value = 35400.1047483396
$ws.Range("B10").Value = 17204.4565791394
$ws.Range("B11").Value = 3038.4757041105
$ws.Range("B12").Value = 35890.176252811
$ws.Range("B13").Value = 643.8459715212
$ws.Range("B14").Value = 1283.5196095865
$ws.Range("B15").Value = 2829.9271633831
$ws.Range("B16").Value = 2605.1478495455
$ws.Range("B17").Value = 4685.2405258054
$ws.Range("B18").Value = 15884.7919834354
$ws.Range("B19").Value = 9130.2411873685
$ws.Range("B20").Value = 7312.0989503456
$ws.Range("B21").Value = 1249.3410186342
$ws.Range("B22").Value = 3492.2337749073
$ws.Range("B23").Value = 1199.7346616256
$ws.Range("B24").Value = 2271.8903002493
$ws.Range("B25").Value = 9056.347303513799
$ws.Range("B26").Value = 1643.6766413736
$ws.Range("B27").Value = 4472.3963103326
$ws.Range("B29").Value = 902.5140603896
$ws.Range("B30").Value = 1293.0703603217
$ws.Range("B31").Value = 2399.8779114579
$ws.Range("B32").Value = 876.691909212
$ws.Range("B33").Value = 60.5506520937
$ws.Range("B34").Value = 788.7682165486
$ws.Range("B35").Value = 537.5141794198

$ws = $wb.Worksheets.Item("steel_sec")
$ws.Range("B2").Value = 3366.0901045362
$ws.Range("B3").Value = 713.2068862166
$ws.Range("B4").Value = 101.8981843381
$ws.Range("B6").Value = 14747.0627807765
$ws.Range("B8").Value = 1676.8916417964
$ws.Range("B9").Value = 13255.9586203289
$ws.Range("B10").Value = 7080.4469558994
$ws.Range("B11").Value = 62.5701237918
$ws.Range("B12").Value = 19485.2690075771
$ws.Range("B14").Value = 3471.2025589502
$ws.Range("B15").Value = 207.1638268035
$ws.Range("B17").Value = 679.8686303652
$ws.Range("B18").Value = 4235.7947524659
$ws.Range("B19").Value = 1906.3681190685
$ws.Range("B20").Value = 1172.5180019987
$ws.Range("B21").Value = 607.5826473323
$ws.Range("B22").Value = 393.8008756541
$ws.Range("B23").Value = 1527.3693102589
$ws.Range("B24").Value = 1549.7034430264
$ws.Range("B25").Value = 1480.3344919248
$ws.Range("B26").Value = 713.115127061
$ws.Range("B27").Value = 1534.5434956711
$ws.Range("B28").Value = 104.2800597572
$ws.Range("B29").Value = 243.0278750417
$ws.Range("B31").Value = 216.1178056209

$ws = $wb.Worksheets.Item("alu_prim")
$ws.Range("B6").Value = 643.5408371324
$ws.Range("B8").Value = 174.877150251
$ws.Range("B9").Value = 468.5192795935
$ws.Range("B10").Value = 395.8014321344
$ws.Range("B12").Value = 0
$ws.Range("B15").Value = 0
$ws.Range("B16").Value = 216.7381705494
$ws.Range("B18").Value = 0
$ws.Range("B20").Value = 255.2634130616
$ws.Range("B21").Value = 72.8971139882
$ws.Range("B22").Value = 162.4320213683
$ws.Range("B24").Value = 135.1309119757
$ws.Range("B25").Value = 282.2007898429
$ws.Range("B26").Value = 1574.3293736099
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 77.74596221349999
$ws.Range("B32").Value = 139.9555562382
$ws.Range("B33").Value = 768.6608934091

$ws = $wb.Worksheets.Item("copper_prim")
$ws.Range("B3").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B26").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 0

$ws = $wb.Worksheets.Item("copper_sec")
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0

$ws = $wb.Worksheets.Item("chlorine")
$ws.Range("B2").Value = 916.8196993567
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 120.5942661158
$ws.Range("B6").Value = 4442.5087201082
$ws.Range("B7").Value = 8.094174998
$ws.Range("B8").Value = 20.1772296278
$ws.Range("B9").Value = 435.4868366428
$ws.Range("B10").Value = 1054.4386657702
$ws.Range("B12").Value = 139.8886383467
$ws.Range("B15").Value = 331.0261540441
$ws.Range("B16").Value = 728.0506172216
$ws.Range("B17").Value = 61.4626436925
$ws.Range("B18").Value = 294.4938873837
$ws.Range("B19").Value = 114.9340091543
$ws.Range("B20").Value = 153.7725459084
$ws.Range("B21").Value = 13.1856323864
$ws.Range("B22").Value = 60.9943306868
$ws.Range("B23").Value = 78.62955662589999
$ws.Range("B24").Value = 59.2409360557
$ws.Range("B25").Value = 420.9457733579
$ws.Range("B26").Value = 291.4030808861
$ws.Range("B27").Value = 16.0172341842
